# The deck currently carries the "Integral" (Red Violet) design on the
# slide master's theme (ppt/theme/theme1.xml) while the notes master's
# theme (ppt/theme/theme2.xml) is the stock "Office Theme" palette.
#
# The edit switches the presentation's applied design back to the
# built-in "Office Theme" colour palette. (Font scheme and format
# scheme are already identical between the two themes, so only the
# 12 theme colours need to change.)
#
# PowerPoint's theme colours are addressed positionally through
# ThemeColorScheme.Colors(1..12) in this fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# and each ThemeColor's .RGB is a standard VBA RGB() packed integer
# (R + G*256 + B*65536), same as the real PowerPoint object model.

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette, in ThemeColorScheme.Colors() slot order.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation

# ThemeColorScheme is exposed per-slide but backs the single shared
# design theme used by the slide master, so any slide works as the
# entry point.
$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgb($officeThemeColors[$i - 1])
}
